$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.969.52'
$ws.Range("E2").Value = '  +0.32%  '
$ws.Range("D3").Value = '1.642.66'
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '''215.69'
$ws.Range("E5").Value = '  +0.17%  '
$ws.Range("E6").Value = '  +1.07%  '
$ws.Range("E7").Value = '  +0.19%  '
$ws.Range("D8").Value = '''0.2566'
$ws.Range("E8").Value = '  +0.23%  '
$ws.Range("E9").Value = '  +0.13%  '
$ws.Range("E10").Value = '  -0.69%  '
$ws.Range("D11").Value = '''0.07777'
$ws.Range("E11").Value = '  +0.60%  '
$ws.Range("D12").Value = '''4.302'
$ws.Range("D13").Value = '1.649.42'
$ws.Range("E13").Value = '  +0.64%  '
$ws.Range("D14").Value = '''0.5479'
$ws.Range("E14").Value = '  +0.55%  '
$ws.Range("D15").Value = '0.0₅7857'
$ws.Range("E15").Value = '  -0.43%  '
$ws.Range("E16").Value = '  +0.19%  '
$ws.Range("D17").Value = '26.028.50'
$ws.Range("E17").Value = '  +0.59%  '
$ws.Range("D18").Value = '''1.004'
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("D19").Value = '''198.44'
$ws.Range("E19").Value = '  -2.38%  '
$ws.Range("D20").Value = '''4.457'
$ws.Range("E20").Value = '  +1.96%  '
$ws.Range("D21").Value = '''9.985'
$ws.Range("E21").Value = '  +1.03%  '
$ws.Range("E22").Value = '  +1.70%  '
$ws.Range("D23").Value = '''1.007'
$ws.Range("E23").Value = '  +0.51%  '
$ws.Range("D24").Value = '''1.875'
$ws.Range("E24").Value = '  -2.89%  '
$ws.Range("D25").Value = '''141.34'
$ws.Range("E25").Value = '  +0.46%  '
$ws.Range("D26").Value = '''0.1155'
$ws.Range("E26").Value = '  +1.85%  '
$ws.Range("D27").Value = '''6.886'
$ws.Range("E27").Value = '  +1.99%  '
$ws.Range("E28").Value = '  +0.48%  '
$ws.Range("D29").Value = '''1.243'
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("D30").Value = '''0.05043'
$ws.Range("E30").Value = '  +1.92%  '
$ws.Range("D31").Value = '''3.266'
$ws.Range("E31").Value = '  -0.09%  '
$ws.Range("D32").Value = '''3.196'
$ws.Range("E32").Value = '  +0.50%  '
$ws.Range("D33").Value = '''1.546'
$ws.Range("E33").Value = '  +0.36%  '
$ws.Range("D34").Value = '''2.365'
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("D35").Value = '''0.9009'
$ws.Range("E35").Value = '  +1.09%  '
$ws.Range("D36").Value = '''2.582'
$ws.Range("E36").Value = '  -1.91%  '
$ws.Range("D37").Value = '1.136.77'
$ws.Range("E37").Value = '  -1.91%  '
$ws.Range("D38").Value = '''0.5507'
$ws.Range("E38").Value = '  -1.85%  '
$ws.Range("E39").Value = '  +15.40%  '
$ws.Range("D40").Value = '''0.01562'
$ws.Range("E40").Value = '  -0.24%  '
$ws.Range("E41").Value = '  +0.49%  '
$ws.Range("D42").Value = '''2.546'
$ws.Range("E42").Value = '  -0.39%  '
$ws.Range("D43").Value = '''5.629'
$ws.Range("E43").Value = '  -0.15%  '
$ws.Range("D44").Value = '''0.8206'
$ws.Range("E44").Value = '  +2.00%  '
$ws.Range("D45").Value = '''100.25'
$ws.Range("E45").Value = '  +0.49%  '
$ws.Range("D46").Value = '1.779.39'
$ws.Range("E46").Value = '  +0.24%  '
$ws.Range("D47").Value = '''0.4536'
$ws.Range("E47").Value = '  -0.27%  '
$ws.Range("D48").Value = '''1.004'
$ws.Range("E48").Value = '  +0.39%  '
$ws.Range("D49").Value = '''55.03'
$ws.Range("E49").Value = '  +0.47%  '
$ws.Range("D50").Value = '''0.05074'
$ws.Range("E50").Value = '  +0.27%  '
$ws.Range("D51").Value = '''1.007'
$ws.Range("E51").Value = '  +0.53%  '
